$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the slightly imprecise timestamp value already stored in A79
$ws.Range("A79").Value = 44392.77451014351

# Append new row 80 with the newly retrieved data
$ws.Range("A80").Value = 44393.77293260515
$ws.Range("B80").Value = 80726
$ws.Range("C80").Value = 68017
$ws.Range("D80").Value = 3658
$ws.Range("E80").Value = 2241
$ws.Range("F80").Value = 1612
$ws.Range("G80").Value = 21331
$ws.Range("H80").Value = 1568
$ws.Range("I80").Value = 915
$ws.Range("J80").Value = 193

# Match the date-style formatting used by the other rows in column A
$ws.Range("A80").NumberFormat = $ws.Range("A79").NumberFormat
